$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the worksheet (this updates the defined name references automatically)
$ws.Name = "Estimates"

# Update the label in A38 from "Min (P=99%)" to "Max (P=99%)"
$ws.Range("A38").Value = "Max (P=99%)"
